$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 253.95454
$ws.Range("I33").Value = 199.39682
$ws.Range("K33").Value = 199.39682
$ws.Range("M33").Value = 29.60318000000001
$ws.Range("H64").Value = 3375
$ws.Range("I64").Value = 3238.889
$ws.Range("J64").Value = 3620
$ws.Range("K64").Value = 3238.889
$ws.Range("L64").Value = 3620
$ws.Range("M64").Value = -2990.889
$ws.Range("N64").Value = -4116
$ws.Range("H67").Value = 3375
$ws.Range("I67").Value = 3238.889
$ws.Range("J67").Value = 3620
$ws.Range("K67").Value = 3238.889
$ws.Range("L67").Value = 3620
$ws.Range("M67").Value = -2380.889
$ws.Range("N67").Value = -5336
$ws.Range("H100").Value = 2742.8333
$ws.Range("I100").Value = 1962
$ws.Range("J100").Value = 3133.25
$ws.Range("K100").Value = 1962
$ws.Range("L100").Value = 3133.25
$ws.Range("M100").Value = -1421
$ws.Range("N100").Value = -4215.25
$ws.Range("H103").Value = 1149
$ws.Range("I103").Value = 1149
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 3447
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -2861
$ws.Range("H107").Value = 1155.1111
$ws.Range("I107").Value = 1016.6667
$ws.Range("J107").Value = 1432
$ws.Range("K107").Value = 1016.6667
$ws.Range("L107").Value = 1432
$ws.Range("M107").Value = 903.3333
$ws.Range("N107").Value = -5272
$ws.Range("H109").Value = 37921
$ws.Range("J109").Value = 37921
$ws.Range("L109").Value = 37921
$ws.Range("N109").Value = -40695
$ws.Range("H112").Value = 2036
$ws.Range("J112").Value = 2062.5642
$ws.Range("L112").Value = 6187.692599999999
$ws.Range("N112").Value = -8403.692599999998
$ws.Range("H118").Value = 1203.6875
$ws.Range("I118").Value = 1322.8572
$ws.Range("J118").Value = 1111
$ws.Range("K118").Value = 3968.5716
$ws.Range("L118").Value = 3333
$ws.Range("M118").Value = -2311.5716
$ws.Range("N118").Value = -6647
$ws.Range("H123").Value = 35726.668
$ws.Range("J123").Value = 35726.668
$ws.Range("L123").Value = 35726.668
$ws.Range("N123").Value = -45526.668
$ws.Range("H128").Value = 58999
$ws.Range("J128").Value = 58999
$ws.Range("L128").Value = 58999
$ws.Range("N128").Value = -68959
$ws.Range("H129").Value = 898.2929
$ws.Range("J129").Value = 926.04443
$ws.Range("L129").Value = 2778.13329
$ws.Range("N129").Value = -12778.13329
$ws.Range("H141").Value = 3793.077
$ws.Range("I141").Value = 1999.4445
$ws.Range("J141").Value = 7828.75
$ws.Range("K141").Value = 5998.333500000001
$ws.Range("L141").Value = 23486.25
$ws.Range("M141").Value = -818.3335000000006
$ws.Range("N141").Value = -33846.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20594.1
$ws.Range("I32").Value = 18992.654
$ws.Range("J32").Value = 31003.5
$ws.Range("K32").Value = 18992.654
$ws.Range("L32").Value = 31003.5
$ws.Range("M32").Value = -18705.654
$ws.Range("N32").Value = -31577.5
$ws.Range("H63").Value = 4082
$ws.Range("J63").Value = 5676
$ws.Range("L63").Value = 5676
$ws.Range("N63").Value = -7048
$ws.Range("H66").Value = 4082
$ws.Range("J66").Value = 5676
$ws.Range("L66").Value = 28380
$ws.Range("N66").Value = -35244
$ws.Range("H109").Value = 44688
$ws.Range("J109").Value = 44688
$ws.Range("L109").Value = 44688
$ws.Range("N109").Value = -47462
$ws.Range("H132").Value = 16131037
$ws.Range("I132").Value = 20834548
$ws.Range("K132").Value = 62503644
$ws.Range("M132").Value = -62501114

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
$ws.Range("H86").Value = 2975
$ws.Range("I86").Value = 2760
$ws.Range("J86").Value = 3333.3333
$ws.Range("K86").Value = 2760
$ws.Range("L86").Value = 3333.3333
$ws.Range("M86").Value = -1637
$ws.Range("N86").Value = -5579.3333
$ws.Range("H89").Value = 2975
$ws.Range("I89").Value = 2760
$ws.Range("J89").Value = 3333.3333
$ws.Range("K89").Value = 13800
$ws.Range("L89").Value = 16666.6665
$ws.Range("M89").Value = -8184
$ws.Range("N89").Value = -27898.6665
$ws.Range("H134").Value = 2517.8308
$ws.Range("I134").Value = 1526.7059
$ws.Range("K134").Value = 4580.1177
$ws.Range("M134").Value = -2045.1177

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1581.1765
$ws.Range("I58").Value = 1156.6774
$ws.Range("J58").Value = 2239.15
$ws.Range("K58").Value = 1156.6774
$ws.Range("L58").Value = 2239.15
$ws.Range("M58").Value = -953.6774
$ws.Range("N58").Value = -2645.15
$ws.Range("H132").Value = 60719.957
$ws.Range("I132").Value = 1692.75
$ws.Range("J132").Value = 178774.38
$ws.Range("K132").Value = 5078.25
$ws.Range("L132").Value = 536323.14
$ws.Range("M132").Value = -2548.25
$ws.Range("N132").Value = -541383.14
$ws.Range("H136").Value = 1581.1765
$ws.Range("I136").Value = 1156.6774
$ws.Range("J136").Value = 2239.15
$ws.Range("K136").Value = 3470.0322
$ws.Range("L136").Value = 6717.450000000001
$ws.Range("M136").Value = -920.0322000000001
$ws.Range("N136").Value = -11817.45

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 281.42856
$ws.Range("I23").Value = 297.33334
$ws.Range("J23").Value = 275.06668
$ws.Range("K23").Value = 892.0000200000001
$ws.Range("L23").Value = 825.2000400000001
$ws.Range("M23").Value = -657.0000200000001
$ws.Range("N23").Value = -1295.20004
$ws.Range("H131").Value = 867.29
$ws.Range("I131").Value = 640
$ws.Range("J131").Value = 871.9286
$ws.Range("K131").Value = 1920
$ws.Range("L131").Value = 2615.7858
$ws.Range("M131").Value = 3120
$ws.Range("N131").Value = -12695.7858

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 4887.5
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 7460
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 7460
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -11300

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1219.5555
$ws.Range("I93").Value = 469
$ws.Range("J93").Value = 1820
$ws.Range("K93").Value = 469
$ws.Range("L93").Value = 1820
$ws.Range("M93").Value = 779
$ws.Range("N93").Value = -4316
$ws.Range("H122").Value = 127538.5
$ws.Range("I122").Value = 145401.14
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 436203.42
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -433753.42
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 3023.6453
$ws.Range("I132").Value = 2695.2173
$ws.Range("J132").Value = 3967.875
$ws.Range("K132").Value = 8085.651899999999
$ws.Range("L132").Value = 11903.625
$ws.Range("M132").Value = -5555.651899999999
$ws.Range("N132").Value = -16963.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 683.0714
$ws.Range("I113").Value = 621.9167
$ws.Range("K113").Value = 1865.7501
$ws.Range("M113").Value = 304.2499
$ws.Range("H126").Value = 1589.4584
$ws.Range("I126").Value = 1090.2222
$ws.Range("J126").Value = 3087.1667
$ws.Range("K126").Value = 3270.6666
$ws.Range("L126").Value = 9261.500100000001
$ws.Range("M126").Value = -800.6665999999996
$ws.Range("N126").Value = -14201.5001
$ws.Range("H132").Value = 1890.6818
$ws.Range("I132").Value = 1554.3529
$ws.Range("J132").Value = 3034.2
$ws.Range("K132").Value = 4663.0587
$ws.Range("L132").Value = 9102.599999999999
$ws.Range("M132").Value = -2133.0587
$ws.Range("N132").Value = -14162.6

Write-Output "Applied all changes"